$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 15; existing rows 15-44 shift down to 16-45.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with values, copying the unchanged
# columns from the row that is now directly below it (row 16, formerly row 15)
# and setting the new / updated values per the diff.
$ws.Range("A15").Value = 2
$ws.Range("B15").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44560
$ws.Range("D15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 100112032
$ws.Range("G15").Value = "Zapallo italiano"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 600
$ws.Range("K15").Value = 6000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 6500
$ws.Range("N15").Value = "$/caja 60 unidades"
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 108
$ws.Range("Q15").Value = 60
$ws.Range("R15").Value = "Hortaliza"
